$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1 (Aufgabe-Eintrag rows 11-16) ---
# Row 11: rename the first activity entry to reflect that it is part 1 of 3
$ws.Range("A11").Value2 = "Bearbeitung Aufgabe 3, Teil 1/3"

# Row 13: fill in the previously empty activity row with a new entry
# (new shared string "Bearbeitung Aufgabe 3, Teil 2" is created here)
$ws.Range("A13").Value2 = "Bearbeitung Aufgabe 3, Teil 2"
$ws.Range("B13").Value2 = 50
$ws.Range("D13").Value2 = 44139

# --- Block 2 (Aufgabe-Eintrag rows 23-28) ---
$ws.Range("A25").Value2 = "Bearbeitung Aufgabe 3, Teil 2"
$ws.Range("B25").Value2 = 50
$ws.Range("D25").Value2 = 44139

# --- Block 3 (Aufgabe-Eintrag rows 35-40) ---
$ws.Range("A37").Value2 = "Bearbeitung Aufgabe 3, Teil 2"
$ws.Range("B37").Value2 = 50
$ws.Range("D37").Value2 = 44139

# Update the active selection to reflect where the author last clicked
[void]$ws.Range("G35").Select()
